# Updates cryptos list data (prices, volume %, and a couple of row re-orderings)
# as produced by the GitHub Actions scraper run on Tue Jul  9 11:13:11 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.298.57"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "3.072.68"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'513.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'140.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.435"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'0.372"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").Value = "3.609.68"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("E13").Value = "  +2.78%  "
$ws.Range("D14").Value = "'25.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.81%  "
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "57.431.50"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "3.079.07"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "'6.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").Value = "'13.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("D20").Value = "'8.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "'333.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "'0.499"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").Value = "'65.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").Value = "'0.169"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.11%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "0.0₃0901"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("D28").Value = "'6.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.49%  "
$ws.Range("D29").Value = "'7.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").Value = "'1.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "'20.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").Value = "'1.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").Value = "'154.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("D34").Value = "'27.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.40%  "
$ws.Range("D35").Value = "'4.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.33%  "
$ws.Range("D36").Value = "'5.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'1.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "'0.0673"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "3.120.17"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").Value = "2.289.29"
$ws.Range("E44").Value = "  +4.68%  "
$ws.Range("D45").Value = "'0.0253"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.90%  "
$ws.Range("D46").Value = "'1.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "'0.937"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'5.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.44%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'19.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("D50").Value = "'0.0876"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").Value = "'249.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.60%  "
